# Applies the edit described by the diff:
#  - Removes the "FORMATO 2" title row from the report header block
#  - Removes an extra blank spacer row, moving the "Fecha:" label up one row
#    (re-using the same right/vcenter style that "Hora:" already had)
#  - Net effect: one row is lost overall, row 1 becomes an empty spacer row,
#    and the final data/header row ends up at row 11 instead of row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a brand new blank row above row 1 (pushes everything down by one).
$ws.Rows("1:1").Insert()

# 2. Delete the row that now holds "FORMATO 2" (previously row 5, now row 6).
$ws.Rows("6:6").Delete()

# 3. Delete the now-redundant blank spacer row above "Fecha:" (previously row 9 with the
#    "Fecha:" text at column H; after deletion, what used to be the all-blank row 8 remains
#    and slides into row 8, with the following "Hora:" row sliding up into row 9).
$ws.Range("H9").Copy() | Out-Null
$ws.Rows("9:9").Delete()

# After the two deletes above, the surviving blank row (old row 8) is now row 8, and the
# "Hora:" row (old row 10) is now row 9. Re-create the "Fecha:" label on row 8, column H,
# copying the formatting that the "Hora:" cell (row 9, column H) already carries so both
# labels share the same look.
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H8").Value = "Fecha:"

# Update the view selection to match the saved file.
$ws.Range("N13").Select() | Out-Null
